$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "theta_threshold_range" row (original row 5). This shifts the
# former row 6 ("pie_threshold_range") up into row 5 and drops the now-unused
# shared string automatically.
$ws.Rows("5").Delete()

# Update the Min/Max values for the remaining parameter rows.
$ws.Range("B2").Formula = "5.4"
$ws.Range("C2").Formula = "10"

$ws.Range("B3").Formula = "6"
$ws.Range("C3").Formula = "8.9"

$ws.Range("B4").Formula = "0.9"
$ws.Range("C4").Formula = "1.4"

$ws.Range("B5").Formula = "0"
$ws.Range("C5").Formula = "15"

# The old row 6 (now row 5) kept a special Times New Roman style on column B
# and a default style on column C; align both with the rest of the data rows.
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to match the saved view.
[void]$ws.Range("C4").Select()

# Restore the window size recorded in the workbook view.
$win = $excel.ActiveWindow
$win.Width = 26025
$win.Height = 9690
